$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header cells (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J15
$data = @(
    @(7,8),
    @(7,7),
    @(7,8),
    @(9,9),
    @(9,9),
    @(7,8),
    @(6,7),
    @(8,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(4,4),
    @(7,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
